$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Price" column (F) with header, type, Korean label, and per-item values
$ws.Range("F1").Value = "Price"
$ws.Range("F2").Value = "int"
$ws.Range("F3").Value = "아이템 가격"
$ws.Range("F4").Value = 60
$ws.Range("F5").Value = 70
$ws.Range("F6").Value = 100
$ws.Range("F7").Value = 50
$ws.Range("F8").Value = 65

# Update the active selection to F9
[void]$ws.Range("F9").Select()
